# Iran/Canada weekly-deaths prediction update.
# The source workbook tracked, per "day the prediction was made" (column A),
# a block of weekly predictions (column B). A new prediction run made on
# 2021-01-02 is inserted between the existing 2020-12-26 and 2021-01-09
# blocks, pushing the old 2021-01-09 block (rows 48-50) down to rows 51-53,
# and the new block (rows 48-50) is populated with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at row 48 - this shifts the old rows 48:50
# (the "2021-01-09" block) down to 51:53, carrying their existing content
# and formatting with them.
$ws.Rows.Item(48).Resize(3).Insert()

# --- New block: prediction made 2021-01-02 -------------------------------
$ws.Range("A48").Value = "2021-01-02"
$ws.Range("B48").Value = "03 Jan -- 09 Jan 2021"
$ws.Range("C48").Value = 94.56999999999999
$ws.Range("D48").Value = 269.17
$ws.Range("E48").Value = 174.6
$ws.Range("F48").Value = "KNN"
$ws.Range("G48").Value = 0.43
$ws.Range("H48").Value = 115.44
$ws.Range("I48").Value = 142.48
$ws.Range("J48").Value = 221.51
$ws.Range("K48").Value = 215.24

$ws.Range("A49").Value = "2021-01-02"
$ws.Range("B49").Value = "10 Jan -- 16 Jan 2021"
$ws.Range("D49").Value = 246.46
$ws.Range("F49").Value = "KNN"

$ws.Range("A50").Value = "2021-01-02"
$ws.Range("B50").Value = "17 Jan -- 23 Jan 2021"
$ws.Range("D50").Value = 213.33
$ws.Range("F50").Value = "KNN"
